$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「動物園」" (row 49) was removed from the data set.
# Delete that entire row; Excel will automatically shift every row
# below it (50-191) up by one, which also updates the sheet's
# used-range dimension from A1:C191 to A1:C190.
$ws.Rows("49:49").Delete()
